$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sensitivity_variables")

# Update sensitivity baseline values
$ws.Range("B10").Value = 0.01
$ws.Range("B11").Value = 0.011

# Update "Include" flags from Y to N for n_exp1 (row 40) and n_exp3 (row 45)
$ws.Range("C40").Value = "N"
$ws.Range("C45").Value = "N"

# Update selection to match the new active cell
$ws.Range("B12").Select()
